# The presentation currently has its live theme ("Integral" colors) stored
# in ppt/theme/theme2.xml, applied via the single slide master / all slides.
# The target edit swaps the Integral and Office Theme color palettes between
# the two theme parts -- in effect, the presentation's applied theme becomes
# "Office Theme" colors (and the Office Theme part takes on the Integral
# colors). We reach the live theme through Slide.ThemeColorScheme, which maps
# 1:1 onto the <a:clrScheme> children in declaration order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGBColor.RGB uses the standard OLE COLORREF (0x00BBGGRR) byte order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
